$wb = $excel.ActiveWorkbook

# --- 1. Remove stray empty cell B3 on "ODI Batting" ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B3").ClearContents()

# --- 2. Add a new worksheet "ODI Batting Extra" after the last existing sheet ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- Header row: reuse the same header style used on the other sheets (bold, centered, bordered) ---
$wsBatting.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Data rows ---
# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$data = @(
    @("3932", 7,    "0", "0", "3.31%", "NO"),
    @("4209", $null, $null, $null, $null, "NO"),
    @("4210", 9,    "1", "0", "1.44%", "NO"),
    @("4660", $null, $null, $null, $null, "NO"),
    @("4663", $null, $null, $null, $null, "NO"),
    @("4666", $null, $null, $null, $null, "NO")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]

    # A: MATCH_CODE - stored as plain text, default (unstyled) cell
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]
    $cellA.Style = "Normal"

    # B: BATTING_POSITION - numeric when present, blank otherwise
    if ($null -ne $rowData[1]) {
        $ws.Cells.Item($row, 2).Value = $rowData[1]
    }

    # C: NUM_4 - text when present
    if ($null -ne $rowData[2]) {
        $cellC = $ws.Cells.Item($row, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $rowData[2]
        $cellC.Style = "Normal"
    }

    # D: NUM_6 - text when present
    if ($null -ne $rowData[3]) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $rowData[3]
        $cellD.Style = "Normal"
    }

    # E: PERCENT_RUNS_OF_TOTAL - text when present
    if ($null -ne $rowData[4]) {
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $rowData[4]
        $cellE.Style = "Normal"
    }

    # F: MAN_OF_MATCH - text, always present
    $cellF = $ws.Cells.Item($row, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $rowData[5]
    $cellF.Style = "Normal"
}

# --- Restore the originally active sheet/selection (adding a sheet activates it) ---
[void]$wb.Worksheets.Item("Player Info").Activate()
[void]$wb.Worksheets.Item("Player Info").Range("A1").Select()

